$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first column (the numeric row index in the original sheet, A2:A18/etc.)
# is removed entirely, which shifts every other column one place to the
# left: old B:F -> new A:E. Deleting the column (rather than just clearing
# values) carries the per-cell formatting along with it, so the old
# (unstyled) column B becomes the new, still-unstyled, column A, while the
# header row (which was fully styled across B1:F1) becomes the styled
# A1:E1.
$ws.Columns.Item(1).Delete()

# The header that used to read "MODEL_CONDITION" is renamed (no underscore).
# After the column shift it sits in D1.
$ws.Range("D1").Value = "MODELCONDITION"
